$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old 3-column table (A1:C3) entirely before laying out the new data.
$ws.Range("A1:C3").Clear()

# New Manager -> Name table.
$ws.Range("A1").Value = "MANAGER"
$ws.Range("B1").Value = "NAME"

$ws.Range("A2").Value = "Ridhima"
$ws.Range("B2").Value = "Divya"

$ws.Range("A3").Value = "Ridhima"
$ws.Range("B3").Value = "Shubham"

$ws.Range("A4").Value = "John"
$ws.Range("B4").Value = "Greg"

$ws.Range("A5").Value = "John"
$ws.Range("B5").Value = "Leo"

# Column C is no longer used by the table -- drop its custom width and
# give column A the width column B/C used to have.
$ws.Columns("C").ColumnWidth = 8.43
$ws.Columns("A").ColumnWidth = 17.7109375

$ws.Range("B3").Select()
